$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update selection on Sheet1 (no longer the active/tabSelected sheet)
$ws1.Range("A1:D6").Select() | Out-Null

# Add the new "Sheet2" worksheet after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 with the four results tables
$ws2.Range("A1").Value = "Trump vs Politicians, Unbalanced"
$ws2.Range("B2").Value = "Accuracy"
$ws2.Range("C2").Value = "AUC"
$ws2.Range("D2").Value = "F-Score"
$ws2.Range("A3").Value = "B-NB"
$ws2.Range("B3").Value = 0.90817999999999999
$ws2.Range("C3").Value = 0.4995
$ws2.Range("D3").Value = 0.87
$ws2.Range("A4").Value = "M-NB"
$ws2.Range("B4").Value = 0.92181000000000002
$ws2.Range("C4").Value = 0.71250000000000002
$ws2.Range("D4").Value = 0.92
$ws2.Range("A5").Value = "SVM"
$ws2.Range("B5").Value = 0.90908999999999995
$ws2.Range("C5").Value = 0.5
$ws2.Range("D5").Value = 0.87
$ws2.Range("A6").Value = "J48"
$ws2.Range("B6").Value = 0.93606
$ws2.Range("C6").Value = 0.73380000000000001
$ws2.Range("D6").Value = 0.93
$ws2.Range("A7").Value = "SGD"
$ws2.Range("B7").Value = 0.94
$ws2.Range("C7").Value = 0.754
$ws2.Range("D7").Value = 0.94
$ws2.Range("A9").Value = "Trump vs Politicians, BALANCED"
$ws2.Range("B10").Value = "Accuracy"
$ws2.Range("C10").Value = "AUC"
$ws2.Range("D10").Value = "F-Score"
$ws2.Range("A11").Value = "B-NB"
$ws2.Range("B11").Value = 0.88533333333300002
$ws2.Range("C11").Value = 0.88533333333300002
$ws2.Range("D11").Value = 0.89
$ws2.Range("A12").Value = "M-NB"
$ws2.Range("B12").Value = 0.89383333333299997
$ws2.Range("C12").Value = 0.89383333333299997
$ws2.Range("D12").Value = 0.89
$ws2.Range("A13").Value = "SVM"
$ws2.Range("B13").Value = 0.50849999999999995
$ws2.Range("C13").Value = 0.50849999999999995
$ws2.Range("D13").Value = 0.36
$ws2.Range("A14").Value = "J48"
$ws2.Range("B14").Value = 0.84283333333300003
$ws2.Range("C14").Value = 0.84283333333300003
$ws2.Range("D14").Value = 0.84
$ws2.Range("A15").Value = "SGD"
$ws2.Range("B15").Value = 0.87983333333299996
$ws2.Range("C15").Value = 0.87983333333299996
$ws2.Range("D15").Value = 0.88
$ws2.Range("A17").Value = "Trump vs Obama, BALANCED"
$ws2.Range("B18").Value = "Accuracy"
$ws2.Range("C18").Value = "AUC"
$ws2.Range("D18").Value = "F-Score"
$ws2.Range("A19").Value = "B-NB"
$ws2.Range("B19").Value = 0.98283333333300005
$ws2.Range("C19").Value = 0.98283333333300005
$ws2.Range("D19").Value = 0.98
$ws2.Range("A20").Value = "M-NB"
$ws2.Range("B20").Value = 0.98
$ws2.Range("C20").Value = 0.98
$ws2.Range("D20").Value = 0.98
$ws2.Range("A21").Value = "SVM"
$ws2.Range("B21").Value = 0.567166666667
$ws2.Range("C21").Value = 0.567166666667
$ws2.Range("D21").Value = 0.47
$ws2.Range("A22").Value = "J48"
$ws2.Range("B22").Value = 0.95650000000000002
$ws2.Range("C22").Value = 0.95650000000000002
$ws2.Range("D22").Value = 0.96
$ws2.Range("A23").Value = "SGD"
$ws2.Range("B23").Value = 0.97133333333299998
$ws2.Range("C23").Value = 0.97133333333299998
$ws2.Range("D23").Value = 0.97
$ws2.Range("A25").Value = "Trump vs Supporters, BALANCED"
$ws2.Range("B26").Value = "Accuracy"
$ws2.Range("C26").Value = "AUC"
$ws2.Range("D26").Value = "F-Score"
$ws2.Range("A27").Value = "B-NB"
$ws2.Range("B27").Value = 0.83650000000000002
$ws2.Range("C27").Value = 0.83650000000000002
$ws2.Range("D27").Value = 0.84
$ws2.Range("A28").Value = "M-NB"
$ws2.Range("B28").Value = 0.84233333333299998
$ws2.Range("C28").Value = 0.84233333333299998
$ws2.Range("D28").Value = 0.84
$ws2.Range("A29").Value = "SVM"
$ws2.Range("B29").Value = 0.51349999999999996
$ws2.Range("C29").Value = 0.51349999999999996
$ws2.Range("D29").Value = 0.36
$ws2.Range("A30").Value = "J48"
$ws2.Range("B30").Value = 0.82683333333300002
$ws2.Range("C30").Value = 0.82683333333300002
$ws2.Range("D30").Value = 0.83
$ws2.Range("A31").Value = "SGD"
$ws2.Range("B31").Value = 0.83733333333299997
$ws2.Range("C31").Value = 0.83733333333299997
$ws2.Range("D31").Value = 0.84

# Leave Sheet2 as the active sheet/tab with D31 selected
$ws2.Range("D31").Select() | Out-Null
